$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "309.23"
Set-TextValue "E2" "-3.99%"
Set-TextValue "D3" "48.66"
Set-TextValue "E3" "-3.22%"
Set-TextValue "E4" "-3.10%"
Set-TextValue "D5" "0.07759"
Set-TextValue "E5" "-4.80%"
Set-TextValue "D6" "4.492"
Set-TextValue "E6" "-2.28%"
Set-TextValue "D7" "1.337"
Set-TextValue "E7" "14.81%"
Set-TextValue "D8" "1.560"
Set-TextValue "E8" "-6.57%"
Set-TextValue "D9" "0.1225"
Set-TextValue "E9" "-8.63%"
Set-TextValue "D10" "0.1945"
Set-TextValue "E10" "-0.38%"
Set-TextValue "D11" "0.04636"
Set-TextValue "E11" "1.40%"
Set-TextValue "D12" "0.09319"
Set-TextValue "E12" "-2.77%"
Set-TextValue "D13" "0.1048"
Set-TextValue "E13" "0.03%"
Set-TextValue "E14" "-5.33%"
Set-TextValue "D15" "0.04173"
Set-TextValue "E15" "-3.07%"
Set-TextValue "D16" "0.005864"
Set-TextValue "E16" "0.74%"
Set-TextValue "E17" "-1.63%"
Set-TextValue "E18" "-6.43%"
Set-TextValue "E19" "2.90%"
Set-TextValue "D20" "8.357"
Set-TextValue "E20" "2.62%"
Set-TextValue "D21" "0.1338"
Set-TextValue "D22" "0.3035"
Set-TextValue "E22" "-0.53%"
Set-TextValue "E23" "-2.26%"
Set-TextValue "D24" "0.004082"
Set-TextValue "E24" "-5.10%"
Set-TextValue "D25" "0.0001351"
Set-TextValue "E25" "0.17%"
Set-TextValue "E26" "-3.98%"
Set-TextValue "D38" "0.02583"
Set-TextValue "E38" "-6.39%"
Set-TextValue "E39" "6.59%"
Set-TextValue "D40" "0.01077"
Set-TextValue "E40" "73.86%"
Set-TextValue "D41" "0.007921"
Set-TextValue "E41" "1.75%"
Set-TextValue "D42" "0.1421"
Set-TextValue "E42" "-1.69%"
Set-TextValue "D43" "0.008400"
Set-TextValue "E43" "9.36%"
Set-TextValue "D44" "0.007684"
Set-TextValue "E44" "-13.11%"
Set-TextValue "D45" "0.3111"
Set-TextValue "E45" "-10.91%"
Set-TextValue "D46" "0.00006941"
Set-TextValue "E46" "2.64%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "0.15%"
Set-TextValue "E48" "-7.58%"
Set-TextValue "E49" "0.22%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "0.15%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "0.15%"
